$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (F1:H1) - match style of existing headers (style index used by C1/D1/E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style from an existing header cell (E1) to the new header cells
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-set the values since paste-format shouldn't touch them, but ensure text remains correct
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Fill boolean FALSE values for rows 2-21 in columns F, G, H
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}
